# "Creating/Updating slides for video recordings"
# Bump the body-text font size on the "JavaScript Code Dilemma" slide
# (slide 2, "Content Placeholder 2") from 30pt to 32pt.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(2)
$sh = $s.Shapes.Item("Content Placeholder 2")

$sh.TextFrame.TextRange.Font.Size = 32
